# Applies the "added description in README.md" commit:
#  - Repositions several shapes on slide 1 and slide 2 (a:off x/y tweaks)
#  - Merges the three text runs "Stinger" / " " / "Team" on slide 2's
#    title placeholder into a single run "Stinger Team"
#
# NOTE: Shape.Left/Top/Width/Height are exposed as single-precision (f32)
# points, and EMU = round-down(points_f32 * 12700) when the deck is
# serialized back to OOXML. The literals below were chosen so that,
# after the f32 round-trip, they floor to the exact target EMU values
# from the diff.

$p = $ppt.ActivePresentation

# ---- Slide 1 ----
$s1 = $p.Slides.Item(1)

# Shape "Title 1" (text "Stinger"): off x 611560 -> 323528 (y unchanged)
$s1.Shapes.Item(1).Left = 25.474647521972656

# Shape "Subtitle 2" (text "JavaScript UI and DOM Teamwork"):
#   off x 539552 -> 467544, y 4221088 -> 4293096
$s1.Shapes.Item(2).Left = 36.814491271972656
$s1.Shapes.Item(2).Top  = 338.0390625

# Shape "Subtitle 2" (text "Just Shoot Me"):
#   off x 467544 -> 251520, y 2852936 -> 2780928
$s1.Shapes.Item(4).Left = 19.804725646972656
$s1.Shapes.Item(4).Top  = 218.97071838378906

# Picture 2: off x unchanged 7092280, y 2636912 -> 2564904
$s1.Shapes.Item(5).Top = 201.96095275878906

# ---- Slide 2 ----
$s2 = $p.Slides.Item(2)

# Title placeholder "Заглавие 1": merge runs "Stinger" + " " + "Team"
# into a single run "Stinger Team" (keeps first run's rPr / endParaRPr)
$titleShape = $s2.Shapes.Item(1)
$fullRange = $titleShape.TextFrame.TextRange
$fullRange.Characters(1, $fullRange.Length).Text = "Stinger Team"

# Shape "Текстово поле 20" (text "Zlatka Todorova"):
#   off x unchanged 323528, y 1412776 -> 1484784
$s2.Shapes.Item(13).Top = 116.9121322631836

# Shape "Текстово поле 25" (text "Tito Titov"):
#   off x 1403648 -> 1331640, y unchanged 4221088
$s2.Shapes.Item(17).Left = 104.85354614257812
